$d = $word.ActiveDocument

# 1. Clear the "EXTENDED BY USE CASE Handle Break-In" text from the
#    Dependency row in the use-case properties table, leaving the
#    paragraph (and its formatting) empty.
$depTable = $d.Tables.Item(2)
$depCell = $depTable.Rows.Item(6).Cells.Item(2)
$depCell.Range.Find.Execute("EXTENDED BY USE CASE Handle Break-In", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 2)

# 2. Remove the entire Basic Flow row that states the system extends by
#    use case Handle Break-In when the timeout is not cancelled.
$flowTable = $d.Tables.Item(3)
$flowTable.Rows.Item(8).Delete()
